# Fix mistake in clue layout: two rooms were both labeled "O".
# The room occupying M13:S15 (with "O"/"O#"/"O*" labels) is actually a
# different room from the real "O" room (columns B:F, rows 16-24) and
# should be relabeled "F" / "F#" / "F*".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13
$ws.Range("M13").Value = "F"
$ws.Range("N13").Value = "F"
$ws.Range("O13").Value = "F"
$ws.Range("P13").Value = "F#"
$ws.Range("Q13").Value = "F"
$ws.Range("R13").Value = "F"
$ws.Range("S13").Value = "F"

# Row 14
$ws.Range("M14").Value = "F"
$ws.Range("N14").Value = "F"
$ws.Range("O14").Value = "F"
$ws.Range("P14").Value = "F*"
$ws.Range("Q14").Value = "F"
$ws.Range("R14").Value = "F"
$ws.Range("S14").Value = "F"

# Row 15
$ws.Range("M15").Value = "F"
$ws.Range("N15").Value = "F"
$ws.Range("O15").Value = "F"
$ws.Range("P15").Value = "F"
$ws.Range("Q15").Value = "F"
$ws.Range("R15").Value = "F"
$ws.Range("S15").Value = "F"

# A handful of these cells had slightly inconsistent formatting
# (plain / alternate-Arial font variants on the same green fill).
# Normalize them to the same black-on-green look used by the rest of
# the room so the whole block is visually consistent.
$normalizeCells = @("Q13", "R13", "S13", "R14", "P15", "Q15", "R15", "S15")
foreach ($addr in $normalizeCells) {
    $cell = $ws.Range($addr)
    $cell.Font.Color = 0
    $cell.Interior.Color = 5220458
}
